# Insert a new data row before the existing row 267 (shifts all rows
# 267..389 down to 268..390, growing the sheet's used range by one row,
# from A1:R389 to A1:R390).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(267).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(267, 1).Value = 10
$ws.Cells.Item(267, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(267, 3).Value = "La Araucanía"
$ws.Cells.Item(267, 4).Value = 44825
$ws.Cells.Item(267, 5).Value = 9
$ws.Cells.Item(267, 6).Value = 100112009
$ws.Cells.Item(267, 7).Value = "Acelga"
$ws.Cells.Item(267, 8).Value = "Sin especificar"
$ws.Cells.Item(267, 9).Value = "Primera"
$ws.Cells.Item(267, 10).Value = 50
$ws.Cells.Item(267, 11).Value = 12000
$ws.Cells.Item(267, 12).Value = 12000
$ws.Cells.Item(267, 13).Value = 12000
$ws.Cells.Item(267, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(267, 15).Value = "Región Metropolitana"
$ws.Cells.Item(267, 16).Value = 1000
$ws.Cells.Item(267, 17).Value = 12
$ws.Cells.Item(267, 18).Value = "Hortaliza"
